$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the table with three new rows (17-19), copying the style of the last
#     existing data row (A16) onto the new index cells in column A ---
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17

# --- Set B (scheme name) and C:M (averaged intensity data) for rows 10-19 with the
#     results of the new run that now includes the Gaussian-Quadrature and three
#     Spiral schemes ---
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.013973233708769
$ws.Range("D10").Value = 1.017022648091769
$ws.Range("E10").Value = 0.9817715755114074
$ws.Range("F10").Value = 1.013973233708769
$ws.Range("G10").Value = 1.005780869297612
$ws.Range("H10").Value = 0.9669194525160409
$ws.Range("I10").Value = 0.9888742712335467
$ws.Range("J10").Value = 1.017022648091769
$ws.Range("K10").Value = 0.9993971118015881
$ws.Range("L10").Value = 1.006685172755178
$ws.Range("M10").Value = 0.9957236750598574

$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9541557546462408
$ws.Range("D11").Value = 0.9051598906807141
$ws.Range("E11").Value = 1.05106765681478
$ws.Range("F11").Value = 0.9541557546462408
$ws.Range("G11").Value = 0.9152379308467258
$ws.Range("H11").Value = 1.174358939256689
$ws.Range("I11").Value = 1.024240379603907
$ws.Range("J11").Value = 0.9051598906807141
$ws.Range("K11").Value = 0.9781137737477469
$ws.Range("L11").Value = 0.9661347641969938
$ws.Range("M11").Value = 1.004036758641509

$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9537777348723918
$ws.Range("D12").Value = 0.9065957406367611
$ws.Range("E12").Value = 1.050732568105727
$ws.Range("F12").Value = 0.9537777348723918
$ws.Range("G12").Value = 0.9160690745303794
$ws.Range("H12").Value = 1.173420721155153
$ws.Range("I12").Value = 1.023863670131787
$ws.Range("J12").Value = 0.9065957406367611
$ws.Range("K12").Value = 0.9786641543712442
$ws.Range("L12").Value = 0.9662209446218181
$ws.Range("M12").Value = 1.004076584905367

$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9540366251249097
$ws.Range("D13").Value = 0.905543009300234
$ws.Range("E13").Value = 1.050981477112953
$ws.Range("F13").Value = 0.9540366251249097
$ws.Range("G13").Value = 0.9154522302873128
$ws.Range("H13").Value = 1.174123122891513
$ws.Range("I13").Value = 1.024121664582482
$ws.Range("J13").Value = 0.905543009300234
$ws.Range("K13").Value = 0.9782622432065933
$ws.Range("L13").Value = 0.9661494341657515
$ws.Range("M13").Value = 1.004043021549901

$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.8679440000000004
$ws.Range("D14").Value = 1.142196000000001
$ws.Range("E14").Value = 1.007944
$ws.Range("F14").Value = 0.8679440000000004
$ws.Range("G14").Value = 1.035784000000002
$ws.Range("H14").Value = 1.062567999999999
$ws.Range("I14").Value = 0.9729439999999989
$ws.Range("J14").Value = 1.142196000000001
$ws.Range("K14").Value = 1.075070000000001
$ws.Range("L14").Value = 0.9715070000000006
$ws.Range("M14").Value = 1.014896666666667

$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.84
$ws.Range("D15").Value = 1.18
$ws.Range("E15").Value = 1.0025
$ws.Range("F15").Value = 0.84
$ws.Range("G15").Value = 1.05
$ws.Range("H15").Value = 1.07
$ws.Range("I15").Value = 0.97
$ws.Range("J15").Value = 1.18
$ws.Range("K15").Value = 1.09125
$ws.Range("L15").Value = 0.965625
$ws.Range("M15").Value = 1.01875

$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9056080435199987
$ws.Range("D16").Value = 1.102067687424001
$ws.Range("E16").Value = 0.9992009854975988
$ws.Range("F16").Value = 0.9056080435199987
$ws.Range("G16").Value = 1.026958587494399
$ws.Range("H16").Value = 1.036987244748796
$ws.Range("I16").Value = 0.9797662181376038
$ws.Range("J16").Value = 1.102067687424001
$ws.Range("K16").Value = 1.0506343364608
$ws.Range("L16").Value = 0.9781211899903992
$ws.Range("M16").Value = 1.008431461137066

$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.99422052270783
$ws.Range("D17").Value = 0.9946418968358763
$ws.Range("E17").Value = 0.9951877037709357
$ws.Range("F17").Value = 0.99422052270783
$ws.Range("G17").Value = 0.9946777625534013
$ws.Range("H17").Value = 0.995197081771158
$ws.Range("I17").Value = 0.9941698166604225
$ws.Range("J17").Value = 0.9946418968358763
$ws.Range("K17").Value = 0.9949148003034061
$ws.Range("L17").Value = 0.994567661505618
$ws.Range("M17").Value = 0.9946824640499372

$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.000301437703879
$ws.Range("D18").Value = 0.9888377589176028
$ws.Range("E18").Value = 0.9943166234628296
$ws.Range("F18").Value = 1.000301437703879
$ws.Range("G18").Value = 0.9931696558585071
$ws.Range("H18").Value = 0.9922960003577823
$ws.Range("I18").Value = 0.9951530746764646
$ws.Range("J18").Value = 0.9888377589176028
$ws.Range("K18").Value = 0.9915771911902163
$ws.Range("L18").Value = 0.9959393144470478
$ws.Range("M18").Value = 0.9940124251628442

$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9996730286711791
$ws.Range("D19").Value = 0.9831266005669683
$ws.Range("E19").Value = 0.9968160848624334
$ws.Range("F19").Value = 0.9996730286711791
$ws.Range("G19").Value = 0.9880776302934834
$ws.Range("H19").Value = 0.9994708258729468
$ws.Range("I19").Value = 0.9970501141675018
$ws.Range("J19").Value = 0.9831266005669683
$ws.Range("K19").Value = 0.9899713427147009
$ws.Range("L19").Value = 0.99482218569294
$ws.Range("M19").Value = 0.9940357140724188
